$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with new match data
$ws.Range("A2").Value = "23-03-2025"
$ws.Range("B2").Value = "Sunrisers Hyderabad vs Rajasthan Royals"
$ws.Range("C2").Value = "Sunrisers Hyderabad"
$ws.Range("D2").Value = "Sunrisers Hyderabad"

# Add new row 3 with new match data
$ws.Range("A3").Value = "23-03-2025"
$ws.Range("B3").Value = "Chennai Super Kings vs Mumbai Indians"
$ws.Range("C3").Value = "Chennai Super Kings"
$ws.Range("D3").Value = "Chennai Super Kings"
